$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# Convert numeric choice codes in column B to text values prefixed with "a"
# (formDef.json cannot have integers in the value list).
$ws.Range("B2").Value = "a1"
$ws.Range("B3").Value = "a2"
$ws.Range("B4").Value = "a3"
$ws.Range("B5").Value = "a888"
$ws.Range("B6").Value = "a1"
$ws.Range("B7").Value = "a2"
$ws.Range("B8").Value = "a3"
$ws.Range("B9").Value = "a4"
$ws.Range("B10").Value = "a5"
$ws.Range("B11").Value = "a6"
$ws.Range("B12").Value = "a888"
$ws.Range("B13").Value = "a1"
$ws.Range("B14").Value = "a2"
$ws.Range("B15").Value = "a3"
$ws.Range("B16").Value = "a4"
$ws.Range("B17").Value = "a5"
$ws.Range("B18").Value = "a6"
$ws.Range("B19").Value = "a7"
$ws.Range("B20").Value = "a8"
$ws.Range("B21").Value = "a9"
$ws.Range("B22").Value = "a10"
$ws.Range("B23").Value = "a888"
$ws.Range("B24").Value = "a1"
$ws.Range("B25").Value = "a2"
$ws.Range("B26").Value = "a3"
$ws.Range("B27").Value = "a4"
$ws.Range("B28").Value = "a5"
$ws.Range("B29").Value = "a6"
$ws.Range("B30").Value = "a7"
$ws.Range("B31").Value = "a8"
$ws.Range("B32").Value = "a9"
$ws.Range("B33").Value = "a888"
$ws.Range("B34").Value = "a0"
$ws.Range("B35").Value = "a1"
$ws.Range("B36").Value = "a2"
$ws.Range("B37").Value = "a3"
$ws.Range("B38").Value = "a4"
$ws.Range("B39").Value = "a5"
$ws.Range("B40").Value = "a888"
$ws.Range("B41").Value = "a1"
$ws.Range("B42").Value = "a0"
$ws.Range("B43").Value = "a2"
$ws.Range("B44").Value = "a888"
$ws.Range("B45").Value = "a0"
$ws.Range("B46").Value = "a1"
$ws.Range("B47").Value = "a2"
$ws.Range("B48").Value = "a3"
$ws.Range("B49").Value = "a1"
$ws.Range("B50").Value = "a2"
$ws.Range("B51").Value = "a1"
$ws.Range("B52").Value = "a2"
$ws.Range("B53").Value = "a3"
$ws.Range("B54").Value = "a4"
$ws.Range("B55").Value = "a5"
$ws.Range("B56").Value = "a888"
$ws.Range("B57").Value = "a1"
$ws.Range("B58").Value = "a0"
$ws.Range("B59").Value = "a888"
$ws.Range("B60").Value = "a1"
$ws.Range("B61").Value = "a666"
$ws.Range("B62").Value = "a777"
$ws.Range("B63").Value = "a888"
$ws.Range("B64").Value = "a999"
$ws.Range("B65").Value = "a0"
$ws.Range("B66").Value = "a1"
$ws.Range("B67").Value = "a2"
$ws.Range("B68").Value = "a999"
$ws.Range("B69").Value = "a1"
$ws.Range("B70").Value = "a0"
$ws.Range("B71").Value = "a888"
$ws.Range("B72").Value = "a999"
$ws.Range("B73").Value = "a1"
$ws.Range("B74").Value = "a0"
$ws.Range("B75").Value = "a1"
$ws.Range("B76").Value = "a2"
$ws.Range("B77").Value = "a888"
$ws.Range("B78").Value = "a999"
$ws.Range("B79").Value = "a1"
$ws.Range("B80").Value = "a0"
$ws.Range("B81").Value = "a2"
$ws.Range("B82").Value = "a888"
$ws.Range("B83").Value = "a999"
$ws.Range("B84").Value = "a1"
$ws.Range("B85").Value = "a2"
$ws.Range("B86").Value = "a3"
$ws.Range("B87").Value = "a4"
$ws.Range("B88").Value = "a888"
$ws.Range("B89").Value = "a999"
$ws.Range("B90").Value = "a1"
$ws.Range("B91").Value = "a0"
$ws.Range("B92").Value = "a2"
$ws.Range("B93").Value = "a888"
$ws.Range("B94").Value = "a999"
$ws.Range("B95").Value = "a1"
$ws.Range("B96").Value = "a2"
$ws.Range("B97").Value = "a3"
$ws.Range("B98").Value = "a4"
$ws.Range("B99").Value = "a5"
$ws.Range("B100").Value = "a6"
$ws.Range("B101").Value = "a7"
$ws.Range("B102").Value = "a888"
$ws.Range("B103").Value = "a1"
$ws.Range("B104").Value = "a2"
$ws.Range("B105").Value = "a3"
$ws.Range("B106").Value = "a4"
$ws.Range("B107").Value = "a5"
$ws.Range("B108").Value = "a6"
$ws.Range("B109").Value = "a7"
$ws.Range("B110").Value = "a8"
$ws.Range("B111").Value = "a888"
$ws.Range("B112").Value = "a1"
$ws.Range("B113").Value = "a0"
$ws.Range("B114").Value = "a666"
$ws.Range("B115").Value = "a888"
$ws.Range("B116").Value = "a999"
$ws.Range("B119").Value = "a777"
$ws.Range("B120").Value = "a888"
$ws.Range("B121").Value = "a999"
$ws.Range("B122").Value = "a1"
$ws.Range("B123").Value = "a2"
$ws.Range("B124").Value = "a3"
$ws.Range("B125").Value = "a888"
$ws.Range("B126").Value = "a999"
$ws.Range("B127").Value = "a1"
$ws.Range("B128").Value = "a2"
$ws.Range("B129").Value = "a3"
$ws.Range("B130").Value = "a4"
$ws.Range("B131").Value = "a5"
$ws.Range("B132").Value = "a6"
$ws.Range("B133").Value = "a888"
$ws.Range("B134").Value = "a1"
$ws.Range("B135").Value = "a2"
$ws.Range("B136").Value = "a3"
$ws.Range("B137").Value = "a4"
$ws.Range("B138").Value = "a5"
$ws.Range("B139").Value = "a6"
$ws.Range("B140").Value = "a7"
$ws.Range("B141").Value = "a8"
$ws.Range("B142").Value = "a9"
$ws.Range("B143").Value = "a10"
$ws.Range("B144").Value = "a11"
$ws.Range("B145").Value = "a12"
$ws.Range("B146").Value = "a13"
$ws.Range("B147").Value = "a14"
$ws.Range("B148").Value = "a15"
$ws.Range("B149").Value = "a16"
$ws.Range("B150").Value = "a888"
$ws.Range("B151").Value = "a1"
$ws.Range("B152").Value = "a2"
$ws.Range("B153").Value = "a3"
$ws.Range("B154").Value = "a4"
$ws.Range("B155").Value = "a5"
$ws.Range("B156").Value = "a6"
$ws.Range("B157").Value = "a7"
$ws.Range("B158").Value = "a8"
$ws.Range("B159").Value = "a9"
$ws.Range("B160").Value = "a10"
$ws.Range("B161").Value = "a11"
$ws.Range("B162").Value = "a12"
$ws.Range("B163").Value = "a13"
$ws.Range("B164").Value = "a14"
$ws.Range("B165").Value = "a15"
$ws.Range("B166").Value = "a888"
$ws.Range("B167").Value = "a1"
$ws.Range("B168").Value = "a2"
$ws.Range("B169").Value = "a3"
$ws.Range("B170").Value = "a4"
$ws.Range("B171").Value = "a5"
$ws.Range("B172").Value = "a6"
$ws.Range("B173").Value = "a888"
$ws.Range("B174").Value = "a1"
$ws.Range("B175").Value = "a0"
$ws.Range("B176").Value = "a888"
$ws.Range("B177").Value = "a999"
$ws.Range("B178").Value = "a0"
$ws.Range("B179").Value = "a1"
$ws.Range("B180").Value = "a2"
$ws.Range("B181").Value = "a888"
$ws.Range("B182").Value = "a0"
$ws.Range("B183").Value = "a1"
$ws.Range("B184").Value = "a2"
$ws.Range("B185").Value = "a3"
$ws.Range("B186").Value = "a4"
$ws.Range("B187").Value = "a5"
$ws.Range("B188").Value = "a6"
$ws.Range("B189").Value = "a7"
$ws.Range("B190").Value = "a8"
$ws.Range("B191").Value = "a9"
$ws.Range("B192").Value = "a10"
$ws.Range("B193").Value = "a11"
$ws.Range("B194").Value = "a12"
$ws.Range("B195").Value = "a13"
$ws.Range("B196").Value = "a14"
$ws.Range("B197").Value = "a888"
$ws.Range("B198").Value = "a4"
$ws.Range("B199").Value = "a3"
$ws.Range("B200").Value = "a2"
$ws.Range("B201").Value = "a1"
$ws.Range("B202").Value = "a0"
$ws.Range("B203").Value = "a888"
$ws.Range("B204").Value = "a1"
$ws.Range("B205").Value = "a2"
$ws.Range("B206").Value = "a3"
$ws.Range("B207").Value = "a4"
$ws.Range("B208").Value = "a5"
$ws.Range("B209").Value = "a6"
$ws.Range("B210").Value = "a888"
$ws.Range("B211").Value = "a999"
$ws.Range("B212").Value = "a0"
$ws.Range("B213").Value = "a1"
$ws.Range("B214").Value = "a888"
$ws.Range("B215").Value = "a1"
$ws.Range("B216").Value = "a0"
$ws.Range("B217").Value = "a2"
$ws.Range("B218").Value = "a888"
$ws.Range("B219").Value = "a1"
$ws.Range("B220").Value = "a2"
$ws.Range("B221").Value = "a3"
$ws.Range("B222").Value = "a0"
$ws.Range("B223").Value = "a888"
$ws.Range("B224").Value = "a999"
$ws.Range("B225").Value = "a0"
$ws.Range("B226").Value = "a1"
$ws.Range("B227").Value = "a2"
$ws.Range("B228").Value = "a3"
$ws.Range("B229").Value = "a4"
$ws.Range("B230").Value = "a888"
$ws.Range("B231").Value = "a1"
$ws.Range("B232").Value = "a2"
$ws.Range("B233").Value = "a3"
$ws.Range("B234").Value = "a4"
$ws.Range("B235").Value = "a5"
$ws.Range("B236").Value = "a888"
$ws.Range("B237").Value = "a999"
$ws.Range("B238").Value = "a0"
$ws.Range("B239").Value = "a1"
$ws.Range("B240").Value = "a2"
$ws.Range("B241").Value = "a888"
$ws.Range("B242").Value = "a999"
$ws.Range("B243").Value = "a1"
$ws.Range("B244").Value = "a0"
$ws.Range("B245").Value = "a0"
$ws.Range("B246").Value = "a1"
$ws.Range("B247").Value = "a1"
$ws.Range("B248").Value = "a0"
$ws.Range("B249").Value = "a2"
$ws.Range("B250").Value = "a1"
$ws.Range("B251").Value = "a0"
$ws.Range("B252").Value = "a3"
$ws.Range("B253").Value = "a1"
$ws.Range("B254").Value = "a1"
$ws.Range("B255").Value = "a888"

# Row 117 value was already text ("Months"); only its alignment/style needs to match
# the rest of the now-text column (bottom-aligned instead of vertically centered).
$ws.Range("B117").VerticalAlignment = -4107

# The whole column B (rows 2-255) now holds text, so remove the vertical-center
# alignment that was used for the old numeric style, matching plain wrap-text style.
$ws.Range("B2:B255").VerticalAlignment = -4107

# Column D never held data; drop its custom (50.33) width so it reverts to the
# sheet default width used from column E onward.
$ws.Columns.Item(4).ColumnWidth = 20.666666666666668

# Make "choices" the active/selected sheet with B10 selected (was "settings" before).
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
